# Insert two new rows at position 237 (shifts existing rows 237:365 down to 239:367)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("237:238").Insert()

# Populate the new row 237 (Primera) with a new weekly price record
$ws.Range("A237").Value = 3
$ws.Range("B237").Value = "Femacal de La Calera"
$ws.Range("C237").Value = "Coquimbo"
$ws.Range("D237").Value = 44452
$ws.Range("E237").Value = 5
$ws.Range("F237").Value = 100112023
$ws.Range("G237").Value = "Brócoli"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 1500
$ws.Range("K237").Value = 600
$ws.Range("L237").Value = 600
$ws.Range("M237").Value = 600
$ws.Range("N237").Value = "`$/unidad"
$ws.Range("O237").Value = "Provincia de Quillota"
$ws.Range("P237").Value = 600
$ws.Range("Q237").Value = 1
$ws.Range("R237").Value = "Hortaliza"

# Populate the new row 238 (Segunda) with a new weekly price record
$ws.Range("A238").Value = 3
$ws.Range("B238").Value = "Femacal de La Calera"
$ws.Range("C238").Value = "Coquimbo"
$ws.Range("D238").Value = 44452
$ws.Range("E238").Value = 5
$ws.Range("F238").Value = 100112023
$ws.Range("G238").Value = "Brócoli"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Segunda"
$ws.Range("J238").Value = 1400
$ws.Range("K238").Value = 500
$ws.Range("L238").Value = 500
$ws.Range("M238").Value = 500
$ws.Range("N238").Value = "`$/unidad"
$ws.Range("O238").Value = "Provincia de Quillota"
$ws.Range("P238").Value = 500
$ws.Range("Q238").Value = 1
$ws.Range("R238").Value = "Hortaliza"
